$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("attribute")

# Remove the "redd_count" attribute row (row 8): deletes the entire row,
# shifts all subsequent rows up by one, and selects the row that now
# occupies that position - mirroring what Excel does when a row is
# deleted via the UI.
$ws.Rows.Item(8).Delete()
$ws.Range("A8:XFD8").Select() | Out-Null
